$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-18 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-08-19 Monday", 2) | Out-Null
$d.Content.Find.Execute("393×9=3537", $true, $true, $false, $false, $false, $true, 1, $false, "829×8=6632", 2) | Out-Null
$d.Content.Find.Execute("229×2=458", $true, $true, $false, $false, $false, $true, 1, $false, "980×3=2940", 2) | Out-Null
$d.Content.Find.Execute("988×8=7904", $true, $true, $false, $false, $false, $true, 1, $false, "778×3=2334", 2) | Out-Null
$d.Content.Find.Execute("902×7=6314", $true, $true, $false, $false, $false, $true, 1, $false, "272×2=544", 2) | Out-Null
$d.Content.Find.Execute("349×2=698", $true, $true, $false, $false, $false, $true, 1, $false, "910×5=4550", 2) | Out-Null
$d.Content.Find.Execute("912×6=5472", $true, $true, $false, $false, $false, $true, 1, $false, "756×9=6804", 2) | Out-Null
$d.Content.Find.Execute("946×2=1892", $true, $true, $false, $false, $false, $true, 1, $false, "445×2=890", 2) | Out-Null
$d.Content.Find.Execute("164×9=1476", $true, $true, $false, $false, $false, $true, 1, $false, "534×6=3204", 2) | Out-Null
$d.Content.Find.Execute("237×8=1896", $true, $true, $false, $false, $false, $true, 1, $false, "626×7=4382", 2) | Out-Null
$d.Content.Find.Execute("641×6=3846", $true, $true, $false, $false, $false, $true, 1, $false, "880×7=6160", 2) | Out-Null
$d.Content.Find.Execute("992×2=1984", $true, $true, $false, $false, $false, $true, 1, $false, "518×2=1036", 2) | Out-Null
$d.Content.Find.Execute("660×2=1320", $true, $true, $false, $false, $false, $true, 1, $false, "245×9=2205", 2) | Out-Null
$d.Content.Find.Execute("777×7=5439", $true, $true, $false, $false, $false, $true, 1, $false, "239×7=1673", 2) | Out-Null
$d.Content.Find.Execute("859×2=1718", $true, $true, $false, $false, $false, $true, 1, $false, "944×7=6608", 2) | Out-Null
$d.Content.Find.Execute("591×7=4137", $true, $true, $false, $false, $false, $true, 1, $false, "738×6=4428", 2) | Out-Null
$d.Content.Find.Execute("657×2=1314", $true, $true, $false, $false, $false, $true, 1, $false, "402×4=1608", 2) | Out-Null
$d.Content.Find.Execute("772×6=4632", $true, $true, $false, $false, $false, $true, 1, $false, "991×3=2973", 2) | Out-Null
$d.Content.Find.Execute("154×2=308", $true, $true, $false, $false, $false, $true, 1, $false, "925×5=4625", 2) | Out-Null
$d.Content.Find.Execute("937×4=3748", $true, $true, $false, $false, $false, $true, 1, $false, "168×6=1008", 2) | Out-Null
$d.Content.Find.Execute("362×2=724", $true, $true, $false, $false, $false, $true, 1, $false, "869×3=2607", 2) | Out-Null
$d.Content.Find.Execute("386×4=1544", $true, $true, $false, $false, $false, $true, 1, $false, "412×4=1648", 2) | Out-Null
$d.Content.Find.Execute("377×6=2262", $true, $true, $false, $false, $false, $true, 1, $false, "608×9=5472", 2) | Out-Null
$d.Content.Find.Execute("322×8=2576", $true, $true, $false, $false, $false, $true, 1, $false, "941×4=3764", 2) | Out-Null
$d.Content.Find.Execute("421×3=1263", $true, $true, $false, $false, $false, $true, 1, $false, "798×8=6384", 2) | Out-Null
$d.Content.Find.Execute("222×4=888", $true, $true, $false, $false, $false, $true, 1, $false, "289×8=2312", 2) | Out-Null
